$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.249.90'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.688.62'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.83'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5246'
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2694'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06435'
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.00'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07470'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.693.96'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.562'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5848'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008516'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.60'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '26.302.41'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.961'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.73'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.218'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.81'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.677'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1235'
$ws.Range('E26').Value = '  +5.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.85'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06658'
$ws.Range('E28').Value = '  +14.40%  '
$ws.Range('E29').Value = '  +5.39%  '
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.587'
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.562'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.668'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.028'
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6217'
$ws.Range('E35').Value = '  +3.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.707'
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.303'
$ws.Range('E38').Value = '  +5.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01620'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').Value = '1.101.23'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8865'
$ws.Range('E41').Value = '  +3.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.23'
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').Value = '1.837.33'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000113'
$ws.Range('E45').Value = '  +2.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.78'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.175'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05262'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4300'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.048'
$ws.Range('E51').Value = '  +2.89%  '
